$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 numeric values ---
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 3.288126333333333
$ws.Range("H2").Value = 9.864379
$ws.Range("I2").Value = 0.05813306630866938
$ws.Range("J2").Value = 0.05813306630866937
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 4.021407666666667
$ws.Range("N2").Value = 12.064223
$ws.Range("O2").Value = 0.06269882270324605
$ws.Range("P2").Value = 0.06269882270324605
$ws.Range("Q2").Value = 13.22289644583522
$ws.Range("R2").Value = 119.006068012517
$ws.Range("S2").Value = 0.003644874817683308
$ws.Range("T2").Value = 0.003644874817683307

# --- Row 3: M2 / Ccl12 / Ccr5 / FAPs ---
$ws.Range("A3").Value = "M2"
$ws.Range("B3").Value = "Ccl12"
$ws.Range("C3").Value = "Ccr5"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 3.288126333333333
$ws.Range("H3").Value = 9.864379
$ws.Range("I3").Value = 0.05813306630866938
$ws.Range("J3").Value = 0.05813306630866937
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.4010506666666667
$ws.Range("N3").Value = 1.203152
$ws.Range("O3").Value = 0.006252886235031953
$ws.Range("P3").Value = 0.006252886235031953
$ws.Range("Q3").Value = 1.318705258067555
$ws.Range("R3").Value = 11.868347322608
$ws.Range("S3").Value = 0.0003634994501216786
$ws.Range("T3").Value = 0.0003634994501216785

# --- Row 4: M2 / Ccl12 / Ccr5 / ECs ---
$ws.Range("A4").Value = "M2"
$ws.Range("B4").Value = "Ccl12"
$ws.Range("C4").Value = "Ccr5"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 3.288126333333333
$ws.Range("H4").Value = 9.864379
$ws.Range("I4").Value = 0.05813306630866938
$ws.Range("J4").Value = 0.05813306630866937
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 59.71602933333333
$ws.Range("N4").Value = 179.148088
$ws.Range("O4").Value = 0.931048291061722
$ws.Range("P4").Value = 0.931048291061722
$ws.Range("Q4").Value = 196.3538485730391
$ws.Range("R4").Value = 1767.184637157352
$ws.Range("S4").Value = 0.0541246920408644
$ws.Range("T4").Value = 0.05412469204086439

# --- Row 5: ECs / Ccl12 / Ccr5 / M2 ---
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Ccl12"
$ws.Range("C5").Value = "Ccr5"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 53.27393966666666
$ws.Range("H5").Value = 159.821819
$ws.Range("I5").Value = 0.9418669336913307
$ws.Range("J5").Value = 0.9418669336913306
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 4.021407666666667
$ws.Range("N5").Value = 12.064223
$ws.Range("O5").Value = 0.06269882270324605
$ws.Range("P5").Value = 0.06269882270324605
$ws.Range("Q5").Value = 214.2362294090707
$ws.Range("R5").Value = 1928.126064681637
$ws.Range("S5").Value = 0.05905394788556274
$ws.Range("T5").Value = 0.05905394788556274

# --- Row 6: ECs / Ccl12 / Ccr5 / FAPs ---
$ws.Range("A6").Value = "ECs"
$ws.Range("B6").Value = "Ccl12"
$ws.Range("C6").Value = "Ccr5"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 53.27393966666666
$ws.Range("H6").Value = 159.821819
$ws.Range("I6").Value = 0.9418669336913307
$ws.Range("J6").Value = 0.9418669336913306
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.4010506666666667
$ws.Range("N6").Value = 1.203152
$ws.Range("O6").Value = 0.006252886235031953
$ws.Range("P6").Value = 0.006252886235031953
$ws.Range("Q6").Value = 21.36554901927644
$ws.Range("R6").Value = 192.289941173488
$ws.Range("S6").Value = 0.005889386784910275
$ws.Range("T6").Value = 0.005889386784910274

# --- Row 7: ECs / Ccl12 / Ccr5 / ECs ---
$ws.Range("A7").Value = "ECs"
$ws.Range("B7").Value = "Ccl12"
$ws.Range("C7").Value = "Ccr5"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 53.27393966666666
$ws.Range("H7").Value = 159.821819
$ws.Range("I7").Value = 0.9418669336913307
$ws.Range("J7").Value = 0.9418669336913306
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 59.71602933333333
$ws.Range("N7").Value = 179.148088
$ws.Range("O7").Value = 0.931048291061722
$ws.Range("P7").Value = 0.931048291061722
$ws.Range("Q7").Value = 3181.308143836896
$ws.Range("R7").Value = 28631.77329453207
$ws.Range("S7").Value = 0.8769235990208577
$ws.Range("T7").Value = 0.8769235990208576
